$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (England)
$ws.Range("F2").Value = 6.608
$ws.Range("G2").Value = 9.911
$ws.Range("H2").Value = 7.589
$ws.Range("I2").Value = 3.256
$ws.Range("J2").Value = 7.781
$ws.Range("K2").Value = 14.101
$ws.Range("L2").Value = 13.455
$ws.Range("M2").Value = 6.775
$ws.Range("N2").Value = 6.871
$ws.Range("O2").Value = 76.34699999999999

# Row 3 (Northern Ireland)
$ws.Range("C3").Value = 2.562
$ws.Range("O3").Value = 2.562

# Row 4 (Scotland)
$ws.Range("D4").Value = 14.891
$ws.Range("O4").Value = 14.891

# Row 5 (Wales)
$ws.Range("E5").Value = 6.201
$ws.Range("O5").Value = 6.201

# Row 6 (COL_TOT)
$ws.Range("C6").Value = 2.562
$ws.Range("D6").Value = 14.891
$ws.Range("E6").Value = 6.201
$ws.Range("F6").Value = 6.608
$ws.Range("G6").Value = 9.911
$ws.Range("H6").Value = 7.589
$ws.Range("I6").Value = 3.256
$ws.Range("J6").Value = 7.781
$ws.Range("K6").Value = 14.101
$ws.Range("L6").Value = 13.455
$ws.Range("M6").Value = 6.775
$ws.Range("N6").Value = 6.871
$ws.Range("O6").Value = 100.001
